# Applies the "Add files via upload" edit:
#  - rename existing 3rd sheet "EnemyStats" -> "EnemyClass"
#  - add a brand-new "EnemyStats" sheet right after "EnemyClass"
#  - EnemyClass gains a "Creature" column + 3 new monster rows
#  - the new EnemyStats sheet gets its own (harder) stat block
#  - UserStats row 2 is replaced (Fighter/Sword -> Necromancer/Scythe)
#  - the new EnemyStats sheet becomes the active sheet/tab

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Rename EnemyStats -> EnemyClass, then insert a fresh EnemyStats
#    sheet right after it.
# ---------------------------------------------------------------
$enemyClass = $wb.Worksheets.Item("EnemyStats")
$enemyClass.Name = "EnemyClass"

$enemyStats = $wb.Worksheets.Add($null, $enemyClass)
$enemyStats.Name = "EnemyStats"

# ---------------------------------------------------------------
# 2) UserStats: row 2 becomes Necromancer / Scythe instead of
#    Fighter / Sword.
# ---------------------------------------------------------------
$userStats = $wb.Worksheets.Item("UserStats")
$userStats.Range("A2").Value = "Necromancer"
$userStats.Range("B2").Value = 27
$userStats.Range("C2").Value = 10
$userStats.Range("D2").Value = 17
$userStats.Range("E2").Value = 10
$userStats.Range("F2").Value = "Scythe"

# ---------------------------------------------------------------
# 3) EnemyClass: add a leading "Creature" column and three rows of
#    monster data (Rat / Evil Turtle / Flora).
# ---------------------------------------------------------------
$enemyClass.Range("A1").Value = "Creature"
$enemyClass.Range("B1").Value = "HP"
$enemyClass.Range("C1").Value = "MP"
$enemyClass.Range("D1").Value = "Attack"
$enemyClass.Range("E1").Value = "Defense"

$enemyClass.Range("A2").Value = "Rat"
$enemyClass.Range("B2").Value = 20
$enemyClass.Range("C2").Value = 0
$enemyClass.Range("D2").Value = 8
$enemyClass.Range("E2").Value = 10

$enemyClass.Range("A3").Value = "Evil Turtle"
$enemyClass.Range("B3").Value = 18
$enemyClass.Range("C3").Value = 0
$enemyClass.Range("D3").Value = 15
$enemyClass.Range("E3").Value = 12

$enemyClass.Range("A4").Value = "Flora"
$enemyClass.Range("B4").Value = 22
$enemyClass.Range("C4").Value = 10
$enemyClass.Range("D4").Value = 10
$enemyClass.Range("E4").Value = 7

$enemyClass.Range("A2").Select() | Out-Null

# ---------------------------------------------------------------
# 4) New EnemyStats sheet: same creature roster, tougher numbers.
# ---------------------------------------------------------------
$enemyStats.Range("A1").Value = "Creature"
$enemyStats.Range("B1").Value = "HP"
$enemyStats.Range("C1").Value = "MP"
$enemyStats.Range("D1").Value = "Attack"
$enemyStats.Range("E1").Value = "Defense"

$enemyStats.Range("A2").Value = "Rat"
$enemyStats.Range("B2").Value = 15
$enemyStats.Range("C2").Value = 0
$enemyStats.Range("D2").Value = 8
$enemyStats.Range("E2").Value = 31

$enemyStats.Range("A3").Value = "Evil Turtle"
$enemyStats.Range("B3").Value = 9
$enemyStats.Range("C3").Value = 0
$enemyStats.Range("D3").Value = 15
$enemyStats.Range("E3").Value = 28

$enemyStats.Range("A4").Value = "Flora"
$enemyStats.Range("B4").Value = 15
$enemyStats.Range("C4").Value = 10
$enemyStats.Range("D4").Value = 10
$enemyStats.Range("E4").Value = 32

$enemyStats.Range("D7").Select() | Out-Null

# ---------------------------------------------------------------
# 5) Make the new EnemyStats sheet the active tab.
# ---------------------------------------------------------------
$enemyStats.Activate() | Out-Null
